$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell values derived from the commit diff: each tuple is
# (CellRef, Value). Dates in column D are Excel serial numbers;
# the cell already carries the date number format (style s="2").
$updates = @(
    @("D2", 45051), @("M2", 50), @("N2", 20000), @("O2", 20000), @("P2", 20000), @("S2", 2857),
    @("D3", 45051), @("M3", 40), @("N3", 15000), @("O3", 15000), @("P3", 15000), @("S3", 2143),
    @("D4", 45027), @("M4", 100), @("N4", 20000), @("O4", 20000), @("P4", 20000), @("S4", 2857),
    @("D5", 45027), @("M5", 80), @("N5", 14000), @("O5", 14000), @("P5", 14000), @("S5", 2000),
    @("D6", 45033), @("M6", 50),
    @("D7", 44644), @("M7", 85), @("N7", 14000), @("O7", 14000), @("P7", 14000), @("S7", 2000),
    @("D8", 44315), @("L8", "Especial"), @("M8", 50), @("N8", 14000), @("O8", 14000), @("P8", 14000), @("S8", 2000),
    @("D9", 44315), @("M9", 80), @("N9", 12000), @("O9", 13000), @("P9", 12500), @("S9", 1786),
    @("D10", 44315), @("N10", 10000), @("O10", 11000), @("P10", 10500), @("S10", 1500),
    @("D11", 45029), @("M11", 120), @("N11", 20000), @("O11", 20000), @("P11", 20000), @("S11", 2857),
    @("D12", 44314), @("M12", 20), @("N12", 13000), @("O12", 13000), @("P12", 13000), @("S12", 1857),
    @("D13", 44314), @("M13", 45), @("N13", 11000), @("O13", 11000), @("P13", 11000), @("S13", 1571),
    @("D14", 45035),
    @("D15", 45035), @("M15", 80),
    @("D16", 44641), @("N16", 13000), @("O16", 13000), @("P16", 13000), @("S16", 1857),
    @("D17", 44699), @("L17", "Segunda"), @("M17", 50), @("N17", 12000), @("O17", 12000), @("P17", 12000), @("S17", 1714),
    @("D18", 44302), @("M18", 340), @("N18", 12000), @("O18", 13000), @("P18", 12500), @("R18", "Provincia de Santiago"), @("S18", 1786),
    @("D19", 44657), @("L19", "Primera"), @("M19", 100), @("N19", 13000), @("O19", 13000), @("P19", 13000), @("S19", 1857),
    @("D20", 45037), @("M20", 50), @("N20", 20000), @("O20", 20000), @("P20", 20000), @("S20", 2857),
    @("D21", 45037), @("L21", "Segunda"), @("M21", 50),
    @("D22", 44659), @("L22", "Primera"), @("N22", 15000), @("O22", 15000), @("P22", 15000), @("S22", 2143),
    @("D23", 44659), @("L23", "Segunda"), @("M23", 20),
    @("D24", 44687), @("M24", 100), @("N24", 15000), @("O24", 15000), @("P24", 15000), @("S24", 2143),
    @("D25", 44687), @("M25", 75), @("N25", 12000), @("O25", 12000), @("P25", 12000), @("S25", 1714),
    @("D26", 44300), @("L26", "Primera"), @("M26", 150), @("O26", 13000), @("P26", 12500), @("R26", "Provincia de Santiago"), @("S26", 1786),
    @("D27", 45049), @("M27", 80),
    @("D28", 45049), @("M28", 50), @("N28", 15000), @("O28", 15000), @("P28", 15000), @("S28", 2143),
    @("D29", 44321), @("M29", 140), @("N29", 11000), @("O29", 12000), @("P29", 11500), @("S29", 1643),
    @("D30", 44321), @("L30", "Segunda"), @("M30", 80), @("N30", 8000), @("O30", 8000), @("P30", 8000), @("S30", 1143),
    @("D31", 44643), @("L31", "Primera"), @("M31", 100),
    @("D32", 45020), @("L32", "Segunda"), @("M32", 120), @("N32", 20000), @("O32", 20000), @("P32", 20000), @("R32", "Región Metropolitana"), @("S32", 2857),
    @("D33", 44335), @("N33", 14000), @("O33", 14000), @("P33", 14000), @("S33", 2000),
    @("D34", 45050), @("L34", "Primera"), @("M34", 100), @("N34", 20000), @("O34", 20000), @("P34", 20000), @("S34", 2857),
    @("D35", 45050), @("M35", 50), @("N35", 15000), @("O35", 15000), @("P35", 15000), @("S35", 2143),
    @("D36", 44312), @("M36", 50), @("N36", 13000), @("O36", 13000), @("P36", 13000), @("R36", "Región Metropolitana"), @("S36", 1857),
    @("L37", "Segunda"), @("M37", 20), @("N37", 11000), @("O37", 11000), @("P37", 11000), @("S37", 1571),
    @("D38", 44306), @("L38", "Primera"), @("M38", 50), @("N38", 12000), @("O38", 12000), @("P38", 12000), @("S38", 1714),
    @("D39", 44306), @("L39", "Segunda"), @("M39", 40), @("N39", 9000), @("O39", 9000), @("P39", 9000), @("S39", 1286),
    @("D40", 44694), @("N40", 15000), @("O40", 15000), @("P40", 15000), @("S40", 2143),
    @("D41", 44694), @("M41", 75), @("N41", 12000), @("O41", 12000), @("P41", 12000), @("S41", 1714),
    @("D42", 44685), @("L42", "Primera"), @("M42", 100), @("N42", 15000), @("O42", 15000), @("P42", 15000), @("S42", 2143),
    @("D43", 44685), @("L43", "Segunda"), @("M43", 70), @("O43", 12000), @("P43", 12000), @("R43", "Región Metropolitana"), @("S43", 1714),
    @("D44", 44623), @("M44", 50),
    @("D45", 44623), @("M45", 30), @("N45", 16000), @("O45", 16000), @("P45", 16000), @("S45", 2286),
    @("D46", 44344), @("L46", "Segunda"), @("M46", 50), @("N46", 12000), @("O46", 12000), @("P46", 12000), @("S46", 1714),
    @("D47", 44307), @("L47", "Primera"), @("M47", 70), @("N47", 14000), @("O47", 14000), @("P47", 14000), @("S47", 2000),
    @("D48", 44307), @("L48", "Segunda"), @("M48", 50), @("N48", 10000), @("O48", 10000), @("P48", 10000), @("S48", 1429),
    @("D49", 44690), @("M49", 100), @("N49", 12000), @("O49", 12000), @("P49", 12000), @("S49", 1714),
    @("D50", 44322), @("L50", "Primera"), @("M50", 100), @("N50", 11000), @("O50", 11000), @("P50", 11000), @("S50", 1571),
    @("D51", 44987), @("M51", 120), @("N51", 18000), @("O51", 18000), @("P51", 18000), @("R51", "Provincia de Santiago"), @("S51", 2571),
    @("D52", 44349), @("M52", 70),
    @("D53", 45043), @("M53", 120), @("N53", 22000), @("O53", 22000), @("P53", 22000), @("S53", 3143),
    @("D54", 45043), @("L54", "Segunda"), @("M54", 80),
    @("D55", 44664), @("L55", "Primera"), @("M55", 80), @("N55", 14000), @("O55", 14000), @("P55", 14000), @("S55", 2000),
    @("D56", 44664), @("L56", "Segunda"), @("M56", 50), @("N56", 12000), @("P56", 12000), @("S56", 1714),
    @("D57", 44316), @("L57", "Primera"), @("M57", 40), @("N57", 13000), @("O57", 13000), @("P57", 13000), @("S57", 1857),
    @("D58", 44316), @("L58", "Segunda"), @("N58", 11000), @("O58", 11000), @("P58", 11000), @("S58", 1571),
    @("D59", 45041), @("L59", "Primera"), @("M59", 100), @("N59", 22000), @("O59", 22000), @("P59", 22000), @("S59", 3143),
    @("D62", 45034), @("L62", "Primera"), @("M62", 100), @("N62", 20000), @("O62", 20000), @("P62", 20000), @("S62", 2857),
    @("D63", 45034), @("L63", "Segunda"), @("M63", 70), @("N63", 14000), @("O63", 14000), @("P63", 14000), @("S63", 2000),
    @("D64", 44679), @("L64", "Primera"), @("M64", 150), @("N64", 12000), @("O64", 12000), @("P64", 12000), @("S64", 1714),
    @("D65", 44342), @("L65", "Segunda"), @("M65", 50), @("N65", 12000), @("O65", 12000), @("P65", 12000), @("S65", 1714),
    @("D66", 45030), @("M66", 100),
    @("D67", 45030), @("M67", 80), @("N67", 15000), @("O67", 15000), @("P67", 15000), @("S67", 2143),
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}